$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 62; $r++) {
    $ws.Cells.Item($r, 3).Value = 7318
}

for ($r = 113; $r -le 245; $r++) {
    $ws.Cells.Item($r, 3).Value = 7293
}
